$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (date-like text) to be stored as text, matching the
# original inline-string cells rather than being auto-parsed as dates.
$ws.Range("C2:C64").NumberFormat = "@"

$ws.Range("C2").Value = "01/01/2017"
$ws.Range("D2").Value = 2953
$ws.Range("C3").Value = "01/04/2017"
$ws.Range("D3").Value = 2925
$ws.Range("C4").Value = "01/07/2017"
$ws.Range("D4").Value = 2937
$ws.Range("C5").Value = "01/10/2017"
$ws.Range("D5").Value = 2962
$ws.Range("C6").Value = "01/01/2018"
$ws.Range("D6").Value = 2983
$ws.Range("C7").Value = "01/04/2018"
$ws.Range("D7").Value = 2989
$ws.Range("C8").Value = "01/07/2018"
$ws.Range("D8").Value = 2978
$ws.Range("C9").Value = "01/10/2018"
$ws.Range("D9").Value = 3002
$ws.Range("C10").Value = "01/01/2019"
$ws.Range("D10").Value = 3018
$ws.Range("C11").Value = "01/04/2019"
$ws.Range("D11").Value = 2982
$ws.Range("C12").Value = "01/07/2019"
$ws.Range("D12").Value = 2987
$ws.Range("C13").Value = "01/10/2019"
$ws.Range("D13").Value = 3014
$ws.Range("C14").Value = "01/01/2020"
$ws.Range("D14").Value = 3048
$ws.Range("C15").Value = "01/04/2022"
$ws.Range("D15").Value = 2827
$ws.Range("C16").Value = "01/07/2022"
$ws.Range("D16").Value = 2930
$ws.Range("C17").Value = "01/10/2022"
$ws.Range("D17").Value = 2984
$ws.Range("C18").Value = "01/01/2023"
$ws.Range("D18").Value = 3004
$ws.Range("C19").Value = "01/04/2023"
$ws.Range("D19").Value = 3003
$ws.Range("C20").Value = "01/07/2023"
$ws.Range("D20").Value = 3053
$ws.Range("C21").Value = "01/10/2023"
$ws.Range("D21").Value = 3077
$ws.Range("C22").Value = "01/01/2024"
$ws.Range("D22").Value = 3123
$ws.Range("C23").Value = "01/01/2017"
$ws.Range("D23").Value = 2030
$ws.Range("C24").Value = "01/04/2017"
$ws.Range("D24").Value = 2028
$ws.Range("C25").Value = "01/07/2017"
$ws.Range("D25").Value = 1998
$ws.Range("C26").Value = "01/10/2017"
$ws.Range("D26").Value = 2064
$ws.Range("C27").Value = "01/01/2018"
$ws.Range("D27").Value = 2072
$ws.Range("C28").Value = "01/04/2018"
$ws.Range("D28").Value = 2070
$ws.Range("C29").Value = "01/07/2018"
$ws.Range("D29").Value = 2070
$ws.Range("C30").Value = "01/10/2018"
$ws.Range("D30").Value = 2086
$ws.Range("C31").Value = "01/01/2019"
$ws.Range("D31").Value = 2093
$ws.Range("C32").Value = "01/04/2019"
$ws.Range("D32").Value = 2070
$ws.Range("C33").Value = "01/07/2019"
$ws.Range("D33").Value = 2050
$ws.Range("C34").Value = "01/10/2019"
$ws.Range("D34").Value = 2080
$ws.Range("C35").Value = "01/01/2020"
$ws.Range("D35").Value = 2105
$ws.Range("C36").Value = "01/04/2022"
$ws.Range("D36").Value = 1893
$ws.Range("C37").Value = "01/07/2022"
$ws.Range("D37").Value = 1976
$ws.Range("C38").Value = "01/10/2022"
$ws.Range("D38").Value = 1992
$ws.Range("C39").Value = "01/01/2023"
$ws.Range("D39").Value = 2053
$ws.Range("C40").Value = "01/04/2023"
$ws.Range("D40").Value = 2036
$ws.Range("C41").Value = "01/07/2023"
$ws.Range("D41").Value = 2047
$ws.Range("C42").Value = "01/10/2023"
$ws.Range("D42").Value = 2072
$ws.Range("C43").Value = "01/01/2024"
$ws.Range("D43").Value = 2104
$ws.Range("C44").Value = "01/01/2017"
$ws.Range("D44").Value = 2352
$ws.Range("C45").Value = "01/04/2017"
$ws.Range("D45").Value = 2249
$ws.Range("C46").Value = "01/07/2017"
$ws.Range("D46").Value = 2133
$ws.Range("C47").Value = "01/10/2017"
$ws.Range("D47").Value = 2075
$ws.Range("C48").Value = "01/01/2018"
$ws.Range("D48").Value = 2179
$ws.Range("C49").Value = "01/04/2018"
$ws.Range("D49").Value = 2078
$ws.Range("C50").Value = "01/07/2018"
$ws.Range("D50").Value = 2189
$ws.Range("C51").Value = "01/10/2018"
$ws.Range("D51").Value = 2147
$ws.Range("C52").Value = "01/01/2019"
$ws.Range("D52").Value = 2142
$ws.Range("C53").Value = "01/04/2019"
$ws.Range("D53").Value = 2073
$ws.Range("C54").Value = "01/07/2019"
$ws.Range("D54").Value = 2044
$ws.Range("C55").Value = "01/10/2019"
$ws.Range("D55").Value = 2029
$ws.Range("C56").Value = "01/01/2020"
$ws.Range("D56").Value = 2150
$ws.Range("C57").Value = "01/04/2022"
$ws.Range("D57").Value = 1985
$ws.Range("C58").Value = "01/07/2022"
$ws.Range("D58").Value = 2072
$ws.Range("C59").Value = "01/10/2022"
$ws.Range("D59").Value = 2131
$ws.Range("C60").Value = "01/01/2023"
$ws.Range("D60").Value = 2142
$ws.Range("C61").Value = "01/04/2023"
$ws.Range("D61").Value = 2191
$ws.Range("C62").Value = "01/07/2023"
$ws.Range("D62").Value = 2094
$ws.Range("C63").Value = "01/10/2023"
$ws.Range("D63").Value = 2081
$ws.Range("C64").Value = "01/01/2024"
$ws.Range("D64").Value = 2168
